$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 58
$ws.Range("G2").Value = 228
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 218
$ws.Range("F3").Value = 22
$ws.Range("G3").Value = 56
$ws.Range("J3").Value = 4.886792452
$ws.Range("K3").Value = 4.886792452
$ws.Range("L3").Value = 4.886792452
$ws.Range("M3").Value = 106
$ws.Range("F4").Value = 79
$ws.Range("G4").Value = 116
$ws.Range("J4").Value = 4.962686567
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 4.962686567
$ws.Range("M4").Value = 134
$ws.Range("F5").Value = 50
$ws.Range("G5").Value = 246
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 5
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 216
$ws.Range("F6").Value = 25
$ws.Range("G6").Value = 88
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 232
$ws.Range("G7").Value = 34
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 5
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 232
$ws.Range("F8").Value = 26
$ws.Range("G8").Value = 115
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 5
$ws.Range("M8").Value = 210
$ws.Range("F9").Value = 26
$ws.Range("G9").Value = 118
$ws.Range("J9").Value = 4.935064935
$ws.Range("K9").Value = 4.909090909
$ws.Range("L9").Value = 4.974025974
$ws.Range("M9").Value = 154
$ws.Range("F10").Value = 25
$ws.Range("G10").Value = 62
$ws.Range("J10").Value = 4.885714285
$ws.Range("K10").Value = 4.885714285
$ws.Range("L10").Value = 4.885714285
$ws.Range("M10").Value = 105
$ws.Range("F11").Value = 20
$ws.Range("G11").Value = 54
$ws.Range("J11").Value = 4.8125
$ws.Range("K11").Value = 4.8125
$ws.Range("L11").Value = 4.8125
$ws.Range("M11").Value = 64
$ws.Range("F12").Value = 50
$ws.Range("G12").Value = 86
$ws.Range("M12").Value = 102
$ws.Range("F13").Value = 52
$ws.Range("G13").Value = 197
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = 5
$ws.Range("M13").Value = 199
$ws.Range("F14").Value = 65
$ws.Range("G14").Value = 209
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 5
$ws.Range("L14").Value = 5
$ws.Range("M14").Value = 157
$ws.Range("F15").Value = 13
$ws.Range("G15").Value = 60
$ws.Range("J15").Value = 4.916129032
$ws.Range("K15").Value = 4.93548387
$ws.Range("L15").Value = 4.903225806
$ws.Range("M15").Value = 155
$ws.Range("F16").Value = 68
$ws.Range("G16").Value = 250
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = 5
$ws.Range("L16").Value = 5
$ws.Range("M16").Value = 164
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 139
$ws.Range("J17").Value = 4.95
$ws.Range("K17").Value = 4.835714285
$ws.Range("L17").Value = 4.892857142
$ws.Range("M17").Value = 140
$ws.Range("F18").Value = 39
$ws.Range("G18").Value = 116
$ws.Range("J18").Value = 4.862068965
$ws.Range("K18").Value = 4.919540229
$ws.Range("L18").Value = 4.804597701
$ws.Range("M18").Value = 87
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 61
$ws.Range("J19").Value = 4.869158878
$ws.Range("K19").Value = 4.887850467
$ws.Range("L19").Value = 4.887850467
$ws.Range("M19").Value = 107
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 75
$ws.Range("J20").Value = 4.876404494
$ws.Range("K20").Value = 4.865168539
$ws.Range("L20").Value = 4.865168539
$ws.Range("M20").Value = 89
$ws.Range("F21").Value = 33
$ws.Range("G21").Value = 93
$ws.Range("J21").Value = 4.857142857
$ws.Range("K21").Value = 4.857142857
$ws.Range("L21").Value = 4.857142857
$ws.Range("M21").Value = 84
$ws.Range("F22").Value = 20
$ws.Range("G22").Value = 100
$ws.Range("J22").Value = 4.905063291
$ws.Range("K22").Value = 4.911392405
$ws.Range("L22").Value = 4.879746835
$ws.Range("M22").Value = 158
$ws.Range("F23").Value = 31
$ws.Range("G23").Value = 125
$ws.Range("J23").Value = 4.892857142
$ws.Range("K23").Value = 4.921428571
$ws.Range("L23").Value = 4.885714285
$ws.Range("M23").Value = 140
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 185
$ws.Range("J24").Value = 5
$ws.Range("K24").Value = 5
$ws.Range("L24").Value = 5
$ws.Range("M24").Value = 228
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 82
$ws.Range("J25").Value = 4.87755102
$ws.Range("K25").Value = 4.87755102
$ws.Range("L25").Value = 4.87755102
$ws.Range("M25").Value = 98
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 123
$ws.Range("J26").Value = 4.918238993
$ws.Range("K26").Value = 4.968553459
$ws.Range("L26").Value = 4.918238993
$ws.Range("M26").Value = 159
$ws.Range("F27").Value = 44
$ws.Range("G27").Value = 120
$ws.Range("J27").Value = 4.903448275
$ws.Range("K27").Value = 4.862068965
$ws.Range("L27").Value = 4.910344827
$ws.Range("M27").Value = 145
$ws.Range("F28").Value = 50
$ws.Range("G28").Value = 46
$ws.Range("J28").Value = 4.888888888
$ws.Range("K28").Value = 4.888888888
$ws.Range("L28").Value = 4.888888888
$ws.Range("M28").Value = 108
$ws.Range("F29").Value = 1017
$ws.Range("G29").Value = 3184
